# Update cryptos list values per latest data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.426.60'
$ws.Cells.Item(2, 5).Value = '  +0.85%  '
$ws.Cells.Item(3, 4).Value = '2.553.90'
$ws.Cells.Item(3, 5).Value = '  +4.91%  '
$ws.Cells.Item(4, 5).Value = '  +0.11%  '
$ws.Cells.Item(5, 4).Value = "'572.26"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.20%  '
$ws.Cells.Item(6, 4).Value = "'150.50"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +7.07%  '
$ws.Cells.Item(7, 5).Value = '  +0.07%  '
$ws.Cells.Item(8, 4).Value = "'0.585"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -0.09%  '
$ws.Cells.Item(9, 4).Value = '2.551.75'
$ws.Cells.Item(9, 5).Value = '  +4.90%  '
$ws.Cells.Item(10, 5).Value = '  +1.64%  '
$ws.Cells.Item(11, 5).Value = '  -1.34%  '
$ws.Cells.Item(12, 5).Value = '  +1.31%  '
$ws.Cells.Item(13, 5).Value = '  +2.48%  '
$ws.Cells.Item(14, 4).Value = "'28.09"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  +6.28%  '
$ws.Cells.Item(15, 4).Value = '3.011.13'
$ws.Cells.Item(15, 5).Value = '  +5.09%  '
$ws.Cells.Item(16, 4).Value = '63.338.74'
$ws.Cells.Item(17, 5).Value = '  +0.68%  '
$ws.Cells.Item(18, 4).Value = '2.571.67'
$ws.Cells.Item(18, 5).Value = '  +5.60%  '
$ws.Cells.Item(19, 4).Value = "'11.67"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  +3.61%  '
$ws.Cells.Item(20, 4).Value = "'341.95"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.38%  '
$ws.Cells.Item(21, 5).Value = '  +2.61%  '
$ws.Cells.Item(22, 4).Value = "'6.80"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  -0.12%  '
$ws.Cells.Item(23, 5).Value = '  +0.11%  '
$ws.Cells.Item(24, 4).Value = "'66.08"
$ws.Cells.Item(24, 4).Style = 'Normal'
$ws.Cells.Item(24, 5).Value = '  +0.92%  '
$ws.Cells.Item(25, 2).Value = 'Fetch.AI'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(25, 4).Value = "'1.62"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +5.42%  '
$ws.Cells.Item(26, 2).Value = 'Kaspa'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(26, 4).Value = "'0.169"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -2.04%  '
$ws.Cells.Item(27, 2).Value = 'SuiNetwork'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(27, 4).Value = "'1.54"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +13.40%  '
$ws.Cells.Item(28, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(28, 4).Value = "'1.00"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.04%  '
$ws.Cells.Item(29, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(29, 4).Value = "'8.45"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  +2.71%  '
$ws.Cells.Item(30, 4).Value = "'7.31"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +11.74%  '
$ws.Cells.Item(31, 4).Value = '0.0₃0821'
$ws.Cells.Item(31, 5).Value = '  +3.99%  '
$ws.Cells.Item(32, 4).Value = "'1.88"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +3.52%  '
$ws.Cells.Item(33, 4).Value = "'177.97"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  +2.23%  '
$ws.Cells.Item(34, 4).Value = "'1.59"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value = '  +8.82%  '
$ws.Cells.Item(35, 4).Value = "'424.83"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +12.25%  '
$ws.Cells.Item(36, 4).Value = "'0.406"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +2.14%  '
$ws.Cells.Item(37, 4).Value = "'18.99"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  +2.07%  '
$ws.Cells.Item(38, 4).Value = "'4.46"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -0.01%  '
$ws.Cells.Item(40, 5).Value = '  +3.56%  '
$ws.Cells.Item(41, 5).Value = '  +0.08%  '
$ws.Cells.Item(42, 4).Value = "'39.49"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  -1.17%  '
$ws.Cells.Item(43, 4).Value = "'153.54"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +5.90%  '
$ws.Cells.Item(44, 5).Value = '  +2.93%  '
$ws.Cells.Item(45, 5).Value = '  +1.60%  '
$ws.Cells.Item(46, 4).Value = "'0.612"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  +2.99%  '
$ws.Cells.Item(47, 4).Value = "'0.0967"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  +1.96%  '
$ws.Cells.Item(48, 2).Value = 'VeChain'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(48, 4).Value = "'0.0240"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  +7.04%  '
$ws.Cells.Item(49, 2).Value = 'Hedera'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(49, 4).Value = "'0.0524"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  +0.69%  '
$ws.Cells.Item(50, 4).Value = "'18.60"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  +3.99%  '
$ws.Cells.Item(51, 5).Value = '  +3.58%  '
